$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("108:108").Insert()

$ws.Range("A108").Value = 5
$ws.Range("B108").Value = "Macroferia Regional de Talca"
$ws.Range("C108").Value = "Maule"
$ws.Range("D108").Value = 44694
$ws.Range("E108").Value = 7
$ws.Range("F108").Value = "Fruta"
$ws.Range("G108").Value = 100108
$ws.Range("H108").Value = "Tropicales y subtropicales"
$ws.Range("I108").Value = 100108002
$ws.Range("J108").Value = "Mango"
$ws.Range("K108").Value = "Sin especificar"
$ws.Range("L108").Value = "Primera"
$ws.Range("M108").Value = 244
$ws.Range("N108").Value = 7000
$ws.Range("O108").Value = 7000
$ws.Range("P108").Value = 7000
$ws.Range("Q108").Value = "`$/bandeja 4 kilos"
$ws.Range("R108").Value = "Brasil"
$ws.Range("S108").Value = 1750
$ws.Range("T108").Value = 4
